$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: cells whose new text looks like a plain number (e.g. "229.53")
# would be auto-coerced to a numeric value by normal .Value assignment (matching
# real Excel semantics). The source data keeps these as text (inline strings),
# so we briefly force a text NumberFormat while writing them, then restore the
# cell style to Normal so no visible/structural formatting change remains.

$ws.Range("D2").Value = "38.061.12"
$ws.Range("E2").Value = "  +2.90%  "

$ws.Range("D3").Value = "2.053.88"
$ws.Range("E3").Value = "  +2.75%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.72"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.94%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("E9").Value = "  +3.56%  "

$ws.Range("E10").Value = "  +4.54%  "

$ws.Range("E11").Value = "  +2.71%  "

$ws.Range("D12").Value = "2.357.08"
$ws.Range("E12").Value = "  +2.75%  "

$ws.Range("E13").Value = "  +5.26%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.53%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.753"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.98%  "

$ws.Range("E16").Value = "  +1.62%  "

$ws.Range("D17").Value = "2.055.57"
$ws.Range("E17").Value = "  +2.68%  "

$ws.Range("D18").Value = "37.955.74"
$ws.Range("E18").Value = "  +2.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.36%  "

$ws.Range("D21").Value = "0.0₃0835"
$ws.Range("E21").Value = "  +3.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.49%  "

$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("E24").Value = "  +0.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.75%  "

$ws.Range("E28").Value = "  +7.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.95%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.93%  "

$ws.Range("E31").Value = "  +2.77%  "

$ws.Range("E32").Value = "  +2.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.06"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +11.01%  "

$ws.Range("E34").Value = "  +3.17%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0605"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.59%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.73%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +14.48%  "

$ws.Range("E38").Value = "  +5.65%  "

$ws.Range("E39").Value = "  +0.13%  "

$ws.Range("D40").Value = "1.523.89"
$ws.Range("E40").Value = "  +5.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.72%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.48%  "

$ws.Range("E45").Value = "  +2.01%  "

$ws.Range("E46").Value = "  +1.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +15.20%  "

$ws.Range("E48").Value = "  +3.02%  "

$ws.Range("E49").Value = "  +3.39%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.56%  "

$ws.Range("D51").Value = "2.244.80"
$ws.Range("E51").Value = "  +2.72%  "

# Row 42/43 swap: VeChain moves to rank 42, HuobiToken moves to rank 43 (with updated price)
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0217"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.94%  "

$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.97%  "
